$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Gradients 2 cruise name in cell A2
$ws.Range("A2").Value = "MGL1704"

# Reflect the active selection on cell A2 as it was when the author saved
$ws.Range("A2").Select()
